# Auto-generated Excel COM-interop script
# Applies numeric value updates to ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
# as produced by the scheduled market-price runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2016.0212
$ws.Range("I15").Value = 2016.0212
$ws.Range("K15").Value = 6048.063599999999
$ws.Range("M15").Value = -5879.063599999999
$ws.Range("H18").Value = 500
$ws.Range("I18").Value = 500
$ws.Range("K18").Value = 500
$ws.Range("M18").Value = -216
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("H32").Value = 833.1667
$ws.Range("J32").Value = 833.1667
$ws.Range("L32").Value = 833.1667
$ws.Range("N32").Value = -1485.1667
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("H88").Value = 2600
$ws.Range("I88").Value = 1800
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 1800
$ws.Range("L88").Value = 3000
$ws.Range("M88").Value = -1394
$ws.Range("N88").Value = -3812
$ws.Range("H91").Value = 2600
$ws.Range("I91").Value = 1800
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 1800
$ws.Range("L91").Value = 3000
$ws.Range("M91").Value = -396
$ws.Range("N91").Value = -5808
$ws.Range("H107").Value = 877.913
$ws.Range("J107").Value = 768.6667
$ws.Range("L107").Value = 768.6667
$ws.Range("N107").Value = -4608.6667
$ws.Range("M20").ClearContents()
$ws.Range("M35").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15656.7705
$ws.Range("I32").Value = 17593.492
$ws.Range("J32").Value = 4564.636
$ws.Range("K32").Value = 17593.492
$ws.Range("L32").Value = 4564.636
$ws.Range("M32").Value = -17306.492
$ws.Range("N32").Value = -5138.636

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 27942.666
$ws.Range("I134").Value = 30217.084
$ws.Range("J134").Value = 649.6667
$ws.Range("K134").Value = 90651.25199999999
$ws.Range("L134").Value = 1949.0001
$ws.Range("M134").Value = -88116.25199999999
$ws.Range("N134").Value = -7019.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 25759.5
$ws.Range("I36").Value = 25759.5
$ws.Range("K36").Value = 25759.5
$ws.Range("M36").Value = -25371.5
$ws.Range("H40").Value = 25759.5
$ws.Range("I40").Value = 25759.5
$ws.Range("K40").Value = 25759.5
$ws.Range("M40").Value = -25599.5
$ws.Range("H92").Value = 37000
$ws.Range("J92").Value = 37000
$ws.Range("L92").Value = 37000
$ws.Range("N92").Value = -41992
$ws.Range("H132").Value = 15432.795
$ws.Range("I132").Value = 16677.395
$ws.Range("K132").Value = 50032.185
$ws.Range("M132").Value = -47502.185
$ws.Range("H134").Value = 756.3461
$ws.Range("I134").Value = 756.3461
$ws.Range("K134").Value = 2269.0383
$ws.Range("M134").Value = 265.9616999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 4185.375
$ws.Range("I2").Value = 6686.8667
$ws.Range("J2").Value = 16.222221
$ws.Range("K2").Value = 40121.2002
$ws.Range("L2").Value = 97.333326
$ws.Range("M2").Value = -40008.2002
$ws.Range("N2").Value = -323.333326
$ws.Range("H12").Value = 55.846153
$ws.Range("I12").Value = 30
$ws.Range("J12").Value = 60.545456
$ws.Range("K12").Value = 90
$ws.Range("L12").Value = 181.636368
$ws.Range("M12").Value = 83
$ws.Range("N12").Value = -527.6363679999999
$ws.Range("H13").Value = 138
$ws.Range("I13").Value = 138
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 414
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -246
$ws.Range("H46").Value = 1866.6666
$ws.Range("I46").Value = 800
$ws.Range("J46").Value = 4000
$ws.Range("K46").Value = 2400
$ws.Range("L46").Value = 12000
$ws.Range("M46").Value = -2309
$ws.Range("N46").Value = -12182
$ws.Range("H107").Value = 4593.5557
$ws.Range("I107").Value = 33666.668
$ws.Range("J107").Value = 959.4167
$ws.Range("K107").Value = 101000.004
$ws.Range("L107").Value = 2878.2501
$ws.Range("M107").Value = -99080.00399999999
$ws.Range("N107").Value = -6718.2501
$ws.Range("H121").Value = 3028.1777
$ws.Range("I121").Value = 565
$ws.Range("J121").Value = 3407.1282
$ws.Range("K121").Value = 1695
$ws.Range("L121").Value = 10221.3846
$ws.Range("M121").Value = -385
$ws.Range("N121").Value = -12841.3846
$ws.Range("H131").Value = 107184.2
$ws.Range("I131").Value = 616
$ws.Range("J131").Value = 119870.89
$ws.Range("K131").Value = 1848
$ws.Range("L131").Value = 359612.67
$ws.Range("M131").Value = 3192
$ws.Range("N131").Value = -369692.67
$ws.Range("N13").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 18187.5
$ws.Range("I15").Value = 16500
$ws.Range("J15").Value = 18428.572
$ws.Range("K15").Value = 16500
$ws.Range("L15").Value = 18428.572
$ws.Range("M15").Value = -16212
$ws.Range("N15").Value = -19004.572
$ws.Range("H80").Value = 11338.77
$ws.Range("I80").Value = 19450.834
$ws.Range("J80").Value = 4385.5713
$ws.Range("K80").Value = 19450.834
$ws.Range("L80").Value = 4385.5713
$ws.Range("M80").Value = -18452.834
$ws.Range("N80").Value = -6381.5713
$ws.Range("H81").Value = 18187.5
$ws.Range("I81").Value = 16500
$ws.Range("J81").Value = 18428.572
$ws.Range("K81").Value = 16500
$ws.Range("L81").Value = 18428.572
$ws.Range("M81").Value = -15502
$ws.Range("N81").Value = -20424.572
$ws.Range("H83").Value = 11338.77
$ws.Range("I83").Value = 19450.834
$ws.Range("J83").Value = 4385.5713
$ws.Range("K83").Value = 97254.17
$ws.Range("L83").Value = 21927.8565
$ws.Range("M83").Value = -92262.17
$ws.Range("N83").Value = -31911.8565
$ws.Range("H84").Value = 18187.5
$ws.Range("I84").Value = 16500
$ws.Range("J84").Value = 18428.572
$ws.Range("K84").Value = 49500
$ws.Range("L84").Value = 55285.716
$ws.Range("M84").Value = -44508
$ws.Range("N84").Value = -65269.716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3325.25
$ws.Range("I22").Value = 3600.3333
$ws.Range("K22").Value = 3600.3333
$ws.Range("M22").Value = -3305.3333
$ws.Range("H27").Value = 3325.25
$ws.Range("I27").Value = 3600.3333
$ws.Range("K27").Value = 3600.3333
$ws.Range("M27").Value = -3493.3333
$ws.Range("H46").Value = 1870.0952
$ws.Range("I46").Value = 1737.3334
$ws.Range("K46").Value = 1737.3334
$ws.Range("M46").Value = -1549.3334
$ws.Range("H63").Value = 18000
$ws.Range("J63").Value = 18000
$ws.Range("L63").Value = 18000
$ws.Range("N63").Value = -19498
$ws.Range("H66").Value = 18000
$ws.Range("J66").Value = 18000
$ws.Range("L66").Value = 54000
$ws.Range("N66").Value = -61488

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 3158
$ws.Range("J15").Value = 3158
$ws.Range("L15").Value = 3158
$ws.Range("N15").Value = -3734
